$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, bordered, centered) by copying H1's format.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns, rows 2-15
$data = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 3)
    5  = @(1, 7)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 3)
    12 = @(1, 5)
    13 = @(7, 9)
    14 = @(1, 3)
    15 = @(4, 5)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
